# Trade update - 2026-02-19 22:22 UTC
# Duplicate the existing trade row (row 2) into a new row 3, carrying over
# its formatting, then update the fields that differ for the new trade:
# exit time, exit CE/PE/combined, P&L (USD/INR/%) and the running
# cumulative P&L formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy row 2's formatting (fonts/fills/number formats/borders) down to
#    row 3 so the new row visually matches the existing trade rows.
$ws.Range("A2:W2").Copy()
$ws.Range("A3").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# 2) Fill in row 3's values. Most columns repeat the same market data as
#    row 2; only the exit time/prices, P&L figures and duration differ.
$ws.Range("A3").Value2 = "19-02-2026"
$ws.Range("B3").Value2 = "Thursday"
$ws.Range("C3").Value2 = "17:34"
$ws.Range("D3").Value2 = "03:52"

$ws.Range("E3").Value2 = 66447.89999999999
$ws.Range("F3").Value2 = 66400
$ws.Range("G3").Value2 = 68800
$ws.Range("H3").Value2 = 61000

$ws.Range("I3").Value2 = 15
$ws.Range("J3").Value2 = 13

$ws.Range("K3").Value2 = 65
$ws.Range("L3").Value2 = 10
$ws.Range("M3").Value2 = 75

$ws.Range("N3").Value2 = 18
$ws.Range("O3").Value2 = 1.4
$ws.Range("P3").Value2 = 19.4

$ws.Range("Q3").Value2 = 55.6
$ws.Range("R3").Value2 = 5046.81
$ws.Range("S3").Value2 = 74.09999999999999

$ws.Range("T3").Value2 = "Time Exit (5:15 PM)"
$ws.Range("U3").Value2 = "0h 0m"
$ws.Range("V3").Value2 = "DRY RUN"

# Running cumulative P&L: previous cumulative total plus this trade's P&L.
$ws.Range("W3").Formula = "=W2+R3"
